$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("K") values were regenerated for rows 2-7 (the Strike# -> K rework).
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 1
$ws.Range("G7").Value = 3
